# Insert a new data row at Excel row 143 (shifts existing rows 143..244 down to 144..245)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(143).Insert()

# Populate the newly inserted row. Columns A,B,C,E,F,G,H,I,O,R carry the same
# values the (now shifted) row below it originally had; D,J,K,L,M,N,P,Q are new.
$ws.Cells.Item(143, 1).Value  = 10
$ws.Cells.Item(143, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(143, 3).Value  = "La Araucanía"
$ws.Cells.Item(143, 4).Value  = 45236
$ws.Cells.Item(143, 5).Value  = 9
$ws.Cells.Item(143, 6).Value  = 100114002
$ws.Cells.Item(143, 7).Value  = "Camote"
$ws.Cells.Item(143, 8).Value  = "Sin especificar"
$ws.Cells.Item(143, 9).Value  = "Primera"
$ws.Cells.Item(143, 10).Value = 100
$ws.Cells.Item(143, 11).Value = 32000
$ws.Cells.Item(143, 12).Value = 32000
$ws.Cells.Item(143, 13).Value = 32000
$ws.Cells.Item(143, 14).Value = "`$/caja 18 kilos"
$ws.Cells.Item(143, 15).Value = "Perú"
$ws.Cells.Item(143, 16).Value = 1778
$ws.Cells.Item(143, 17).Value = 18
$ws.Cells.Item(143, 18).Value = "Hortaliza"

# Ensure the date column keeps the same date/time number format as the rest of column D
$ws.Cells.Item(143, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
